$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.8803455971160747
$ws.Range("D2").Value = 0.9402674389889207
$ws.Range("E2").Value = 0.9003264978541532
$ws.Range("F2").Value = 0.5541011900173972
$ws.Range("G2").Value = 25.8844033916166
$ws.Range("H2").Value = 34.59109753736143
$ws.Range("I2").Value = 1.065387804516424
$ws.Range("J2").Value = -0.967385085186681
$ws.Range("C3").Value = 0.7334600689828357
$ws.Range("D3").Value = 0.8570054674833605
$ws.Range("E3").Value = 1.343744905684221
$ws.Range("F3").Value = 1.034531293126445
$ws.Range("G3").Value = 48.32731962134009
$ws.Range("H3").Value = 51.62750536459847
$ws.Range("I3").Value = 1.032986868166707
$ws.Range("J3").Value = -0.5015754546980471
$ws.Range("C4").Value = 0.6998537242161332
$ws.Range("D4").Value = 0.8479258243899839
$ws.Range("E4").Value = 1.230787330732287
$ws.Range("F4").Value = 0.9801530363346175
$ws.Range("G4").Value = 52.73117591895465
$ws.Range("H4").Value = 54.78560721429185
$ws.Range("I4").Value = 0.8690601118660353
$ws.Range("J4").Value = 1.750219356180542
$ws.Range("C5").Value = 0.7147216923100759
$ws.Range("D5").Value = 0.845952123766188
$ws.Range("E5").Value = 0.2986650167366291
$ws.Range("F5").Value = 0.07203614282667244
$ws.Range("G5").Value = 17.05214979678195
$ws.Range("H5").Value = 53.41145080316804
$ws.Range("I5").Value = 0.9823367530696279
$ws.Range("J5").Value = 0.03422265787892975
$ws.Range("C6").Value = 0.3666564107336027
$ws.Range("D6").Value = 0.6365098407590214
$ws.Range("E6").Value = 0.4450098977000913
$ws.Range("F6").Value = 0.2464597688206375
$ws.Range("G6").Value = 58.34111505556153
$ws.Range("H6").Value = 79.58288693346059
$ws.Range("I6").Value = 0.7647380427077389
$ws.Range("J6").Value = 0.2699175300701723
$ws.Range("C7").Value = 0.7059895322335761
$ws.Range("D7").Value = 0.8440246840543375
$ws.Range("E7").Value = 0.2784799137574325
$ws.Range("F7").Value = 0.1777348932259392
$ws.Range("G7").Value = 42.27788918134943
$ws.Range("H7").Value = 54.22273211176508
$ws.Range("I7").Value = 0.9380516873956877
$ws.Range("J7").Value = 0.04004302719950181
$ws.Range("C8").Value = 0.9089349918444056
$ws.Range("D8").Value = 0.9539642675183697
$ws.Range("E8").Value = 0.09449594347390443
$ws.Range("F8").Value = 0.02475677424468299
$ws.Range("G8").Value = 9.956344239826164
$ws.Range("H8").Value = 30.1769793312045
$ws.Range("I8").Value = 1.030680300904954
$ws.Range("J8").Value = -0.09563584612164799
$ws.Range("C9").Value = 0.3575627036804937
$ws.Range("D9").Value = 0.6301498879048755
$ws.Range("E9").Value = 0.2509878931074921
$ws.Range("F9").Value = 0.1605333443212644
$ws.Range("G9").Value = 64.56112667329117
$ws.Range("H9").Value = 80.15218626584719
$ws.Range("I9").Value = 0.7605747739649659
$ws.Range("J9").Value = 0.7935840813035195
$ws.Range("C10").Value = 0.6527787757704895
$ws.Range("D10").Value = 0.8228352472769136
$ws.Range("E10").Value = 0.1454711414789391
$ws.Range("F10").Value = 0.1095819739374891
$ws.Range("G10").Value = 51.7630552866941
$ws.Range("H10").Value = 58.92548041632843
$ws.Range("I10").Value = 0.8407873958158331
$ws.Range("J10").Value = 0.5175436308950983
$ws.Range("C11").Value = 0.3066734816357521
$ws.Range("D11").Value = 0.6841457226781594
$ws.Range("E11").Value = 84.17372061022785
$ws.Range("F11").Value = 65.82941403978124
$ws.Range("G11").Value = 81.39455676364473
$ws.Range("H11").Value = 83.266230751983
$ws.Range("I11").Value = 2.420832657636722
$ws.Range("J11").Value = -733.480472768125
$ws.Range("C12").Value = 0.2776019442249618
$ws.Range("D12").Value = 0.6419864272914642
$ws.Range("E12").Value = 85.9203233206614
$ws.Range("F12").Value = 67.94036118541874
$ws.Range("G12").Value = 84.00463023576789
$ws.Range("H12").Value = 84.9940030693365
$ws.Range("I12").Value = 2.331830650016999
$ws.Range("J12").Value = -687.0465232421046
$ws.Range("C13").Value = 0.4417101623820751
$ws.Range("D13").Value = 0.7922728588326087
$ws.Range("E13").Value = 61.16598027552026
$ws.Range("F13").Value = 53.21522896573224
$ws.Range("G13").Value = 79.32008903403404
$ws.Range("H13").Value = 74.71879533410085
$ws.Range("I13").Value = 2.013607310378205
$ws.Range("J13").Value = -533.4319833474813
$ws.Range("C14").Value = 0.7665023519639267
$ws.Range("D14").Value = 0.8819758924143449
$ws.Range("E14").Value = 1.52047269352845
$ws.Range("F14").Value = 0.5866645174814363
$ws.Range("G14").Value = 24.60223374913872
$ws.Range("H14").Value = 48.32159434829042
$ws.Range("I14").Value = 1.132782949118682
$ws.Range("J14").Value = -1.934677823168608
$ws.Range("C15").Value = 0.661785037146844
$ws.Range("D15").Value = 0.8182022218247786
$ws.Range("E15").Value = 1.829927055992258
$ws.Range("F15").Value = 1.047247764198787
$ws.Range("G15").Value = 43.91715114916035
$ws.Range("H15").Value = 58.15625184390376
$ws.Range("I15").Value = 1.111281491793032
$ws.Range("J15").Value = -1.591755002803323
$ws.Range("C16").Value = 0.7544959021431671
$ws.Range("D16").Value = 0.8751431978866244
$ws.Range("E16").Value = 1.188009207205942
$ws.Range("F16").Value = 0.9229067114792806
$ws.Range("G16").Value = 46.73840268465634
$ws.Range("H16").Value = 49.54837008992655
$ws.Range("I16").Value = 0.8992593967004274
$ws.Range("J16").Value = 1.455656934604308
